$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 854.8182  # H2: '899.7' -> '854.8182'
$ws.Cells.Item(2, 9).Value = 486.14285  # I2: '119.4' -> '486.14285'
$ws.Cells.Item(2, 10).Value = 1500  # J2: '1680' -> '1500'
$ws.Cells.Item(2, 11).Value = 486.14285  # K2: '119.4' -> '486.14285'
$ws.Cells.Item(2, 12).Value = 1500  # L2: '1680' -> '1500'
$ws.Cells.Item(2, 13).Value = -373.14285  # M2: '-6.400000000000006' -> '-373.14285'
$ws.Cells.Item(2, 14).Value = -1726  # N2: '-1906' -> '-1726'

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 206.2  # H38: '26.2' -> '206.2'
$ws.Cells.Item(38, 9).Value = 7.75  # I38: '26.2' -> '7.75'
$ws.Cells.Item(38, 10).Value = 1000  # J38: '0' -> '1000'
$ws.Cells.Item(38, 11).Value = 23.25  # K38: '78.59999999999999' -> '23.25'
$ws.Cells.Item(38, 12).Value = 3000  # L38: '0' -> '3000'
$ws.Cells.Item(38, 13).Value = 348.75  # M38: '293.4' -> '348.75'
$ws.Cells.Item(38, 14).Value = -3744  # N38: None -> '-3744'

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(42, 8).Value = 192.84616  # H42: '206.25' -> '192.84616'
$ws.Cells.Item(42, 9).Value = 79.57143000000001  # I42: '110.71429' -> '79.57143000000001'
$ws.Cells.Item(42, 10).Value = 325  # J42: '340' -> '325'
$ws.Cells.Item(42, 11).Value = 238.71429  # K42: '332.14287' -> '238.71429'
$ws.Cells.Item(42, 12).Value = 975  # L42: '1020' -> '975'
$ws.Cells.Item(42, 13).Value = -8.714290000000005  # M42: '-102.14287' -> '-8.714290000000005'
$ws.Cells.Item(42, 14).Value = -1435  # N42: '-1480' -> '-1435'

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 842  # H53: '772.8' -> '842'
$ws.Cells.Item(53, 9).Value = 922.5  # I53: '836.6667' -> '922.5'
$ws.Cells.Item(53, 11).Value = 922.5  # K53: '836.6667' -> '922.5'
$ws.Cells.Item(53, 13).Value = -285.5  # M53: '-199.6667' -> '-285.5'

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 8372  # H116: '7911.8335' -> '8372'
$ws.Cells.Item(116, 9).Value = 8500  # I116: '7494.3335' -> '8500'
$ws.Cells.Item(116, 11).Value = 8500  # K116: '7494.3335' -> '8500'
$ws.Cells.Item(116, 13).Value = -5058  # M116: '-4052.3335' -> '-5058'

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 6000  # H31: '5499.5' -> '6000'
$ws.Cells.Item(31, 9).Value = 6000  # I31: '5499.5' -> '6000'
$ws.Cells.Item(31, 11).Value = 6000  # K31: '5499.5' -> '6000'
$ws.Cells.Item(31, 13).Value = -5706  # M31: '-5205.5' -> '-5706'

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 6367.636  # H61: '6504.5' -> '6367.636'
$ws.Cells.Item(61, 9).Value = 3528.1667  # I61: '3234' -> '3528.1667'
$ws.Cells.Item(61, 11).Value = 3528.1667  # K61: '3234' -> '3528.1667'
$ws.Cells.Item(61, 13).Value = -3316.1667  # M61: '-3022' -> '-3316.1667'

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(94, 8).Value = 65000  # H94: '0' -> '65000'
$ws.Cells.Item(94, 10).Value = 65000  # J94: '0' -> '65000'
$ws.Cells.Item(94, 12).Value = 65000  # L94: '0' -> '65000'
$ws.Cells.Item(94, 14).Value = -66802  # N94: None -> '-66802'

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 2829.375  # H102: '3337.2' -> '2829.375'
$ws.Cells.Item(102, 9).Value = 2519.2856  # I102: '2921.5' -> '2519.2856'
$ws.Cells.Item(102, 11).Value = 2519.2856  # K102: '2921.5' -> '2519.2856'
$ws.Cells.Item(102, 13).Value = -897.2856000000002  # M102: '-1299.5' -> '-897.2856000000002'

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 9443.799999999999  # H132: '10270.889' -> '9443.799999999999'
$ws.Cells.Item(132, 9).Value = 3487.6  # I132: '3859.5' -> '3487.6'
$ws.Cells.Item(132, 11).Value = 10462.8  # K132: '11578.5' -> '10462.8'
$ws.Cells.Item(132, 13).Value = -7932.799999999999  # M132: '-9048.5' -> '-7932.799999999999'

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 6367.636  # H136: '6504.5' -> '6367.636'
$ws.Cells.Item(136, 9).Value = 3528.1667  # I136: '3234' -> '3528.1667'
$ws.Cells.Item(136, 11).Value = 10584.5001  # K136: '9702' -> '10584.5001'
$ws.Cells.Item(136, 13).Value = -8034.500100000001  # M136: '-7152' -> '-8034.500100000001'

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2624.75  # H99: '4252.5' -> '2624.75'
$ws.Cells.Item(99, 9).Value = 2642.5715  # I99: '4252.5' -> '2642.5715'
$ws.Cells.Item(99, 10).Value = 2500  # J99: '0' -> '2500'
$ws.Cells.Item(99, 11).Value = 2642.5715  # K99: '4252.5' -> '2642.5715'
$ws.Cells.Item(99, 12).Value = 2500  # L99: '0' -> '2500'
$ws.Cells.Item(99, 13).Value = -1144.5715  # M99: '-2754.5' -> '-1144.5715'
$ws.Cells.Item(99, 14).Value = -5496  # N99: None -> '-5496'

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2579.5  # H105: '2664.7144' -> '2579.5'
$ws.Cells.Item(105, 9).Value = 2393.3333  # I105: '2475.4' -> '2393.3333'
$ws.Cells.Item(105, 11).Value = 2393.3333  # K105: '2475.4' -> '2393.3333'
$ws.Cells.Item(105, 13).Value = -646.3332999999998  # M105: '-728.4000000000001' -> '-646.3332999999998'

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1274.75  # H22: '1075' -> '1274.75'
$ws.Cells.Item(22, 9).Value = 1239.2  # I22: '1599' -> '1239.2'
$ws.Cells.Item(22, 10).Value = 1334  # J22: '900.3333' -> '1334'
$ws.Cells.Item(22, 11).Value = 1239.2  # K22: '1599' -> '1239.2'
$ws.Cells.Item(22, 12).Value = 1334  # L22: '900.3333' -> '1334'
$ws.Cells.Item(22, 13).Value = -889.2  # M22: '-1249' -> '-889.2'
$ws.Cells.Item(22, 14).Value = -2034  # N22: '-1600.3333' -> '-2034'

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 715.1  # H107: '452.8889' -> '715.1'
$ws.Cells.Item(107, 9).Value = 857  # I107: '468.2857' -> '857'
$ws.Cells.Item(107, 10).Value = 384  # J107: '399' -> '384'
$ws.Cells.Item(107, 11).Value = 857  # K107: '468.2857' -> '857'
$ws.Cells.Item(107, 12).Value = 384  # L107: '399' -> '384'
$ws.Cells.Item(107, 13).Value = 1063  # M107: '1451.7143' -> '1063'
$ws.Cells.Item(107, 14).Value = -4224  # N107: '-4239' -> '-4224'

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 1383  # H122: '1437.8572' -> '1383'
$ws.Cells.Item(122, 9).Value = 1123.4286  # I122: '1144.1666' -> '1123.4286'
$ws.Cells.Item(122, 11).Value = 3370.2858  # K122: '3432.4998' -> '3370.2858'
$ws.Cells.Item(122, 13).Value = -920.2857999999997  # M122: '-982.4998000000001' -> '-920.2857999999997'

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 10549.1  # H134: '10883.167' -> '10549.1'
$ws.Cells.Item(134, 9).Value = 6337  # I134: '5599.6665' -> '6337'
$ws.Cells.Item(134, 10).Value = 12354.286  # J134: '16166.667' -> '12354.286'
$ws.Cells.Item(134, 11).Value = 19011  # K134: '16798.9995' -> '19011'
$ws.Cells.Item(134, 12).Value = 37062.858  # L134: '48500.001' -> '37062.858'
$ws.Cells.Item(134, 13).Value = -16476  # M134: '-14263.9995' -> '-16476'
$ws.Cells.Item(134, 14).Value = -42132.858  # N134: '-53570.001' -> '-42132.858'

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 53.153847  # H6: '43.42857' -> '53.153847'
$ws.Cells.Item(6, 9).Value = 58.1  # I6: '46' -> '58.1'
$ws.Cells.Item(6, 10).Value = 36.666668  # J6: '37' -> '36.666668'
$ws.Cells.Item(6, 11).Value = 174.3  # K6: '138' -> '174.3'
$ws.Cells.Item(6, 12).Value = 110.000004  # L6: '111' -> '110.000004'
$ws.Cells.Item(6, 13).Value = -61.30000000000001  # M6: '-25' -> '-61.30000000000001'
$ws.Cells.Item(6, 14).Value = -336.000004  # N6: '-337' -> '-336.000004'

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 987.5  # H26: '7243.25' -> '987.5'
$ws.Cells.Item(26, 10).Value = 987.5  # J26: '7243.25' -> '987.5'
$ws.Cells.Item(26, 12).Value = 2962.5  # L26: '21729.75' -> '2962.5'
$ws.Cells.Item(26, 14).Value = -3538.5  # N26: '-22305.75' -> '-3538.5'

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 3251.25  # H55: '4176.25' -> '3251.25'
$ws.Cells.Item(55, 10).Value = 3251.25  # J55: '4176.25' -> '3251.25'
$ws.Cells.Item(55, 12).Value = 9753.75  # L55: '12528.75' -> '9753.75'
$ws.Cells.Item(55, 14).Value = -10107.75  # N55: '-12882.75' -> '-10107.75'

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(93, 8).Value = 0  # H93: '855' -> '0'
$ws.Cells.Item(93, 10).Value = 0  # J93: '855' -> '0'
$ws.Cells.Item(93, 12).Value = 0  # L93: '2565' -> '0'
$ws.Cells.Item(93, 14).ClearContents()  # N93: was '-6309'

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 4999.75  # H113: '4325' -> '4999.75'
$ws.Cells.Item(113, 9).Value = 5333  # I113: '6000' -> '5333'
$ws.Cells.Item(113, 10).Value = 4000  # J113: '2650' -> '4000'
$ws.Cells.Item(113, 11).Value = 5333  # K113: '6000' -> '5333'
$ws.Cells.Item(113, 12).Value = 4000  # L113: '2650' -> '4000'
$ws.Cells.Item(113, 13).Value = -3163  # M113: '-3830' -> '-3163'
$ws.Cells.Item(113, 14).Value = -8340  # N113: '-6990' -> '-8340'

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2027.7778  # H126: '2651' -> '2027.7778'
$ws.Cells.Item(126, 9).Value = 2542.1667  # I126: '2651' -> '2542.1667'
$ws.Cells.Item(126, 10).Value = 999  # J126: '0' -> '999'
$ws.Cells.Item(126, 11).Value = 7626.500100000001  # K126: '7953' -> '7626.500100000001'
$ws.Cells.Item(126, 12).Value = 2997  # L126: '0' -> '2997'
$ws.Cells.Item(126, 13).Value = -5156.500100000001  # M126: '-5483' -> '-5156.500100000001'
$ws.Cells.Item(126, 14).Value = -7937  # N126: None -> '-7937'

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 140837.4  # H132: '163076' -> '140837.4'
$ws.Cells.Item(132, 9).Value = 170547.17  # I132: '205399.3' -> '170547.17'
$ws.Cells.Item(132, 11).Value = 511641.51  # K132: '616197.8999999999' -> '511641.51'
$ws.Cells.Item(132, 13).Value = -509111.51  # M132: '-613667.8999999999' -> '-509111.51'

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2180  # H22: '1998.3334' -> '2180'
$ws.Cells.Item(22, 9).Value = 2100  # I22: '1898' -> '2100'
$ws.Cells.Item(22, 11).Value = 2100  # K22: '1898' -> '2100'
$ws.Cells.Item(22, 13).Value = -1805  # M22: '-1603' -> '-1805'

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 2180  # H27: '1998.3334' -> '2180'
$ws.Cells.Item(27, 9).Value = 2100  # I27: '1898' -> '2100'
$ws.Cells.Item(27, 11).Value = 2100  # K27: '1898' -> '2100'
$ws.Cells.Item(27, 13).Value = -1993  # M27: '-1791' -> '-1993'

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2082.182  # H82: '1900.625' -> '2082.182'
$ws.Cells.Item(82, 9).Value = 1866.8334  # I82: '1540.4' -> '1866.8334'
$ws.Cells.Item(82, 10).Value = 2340.6  # J82: '2501' -> '2340.6'
$ws.Cells.Item(82, 11).Value = 1866.8334  # K82: '1540.4' -> '1866.8334'
$ws.Cells.Item(82, 12).Value = 2340.6  # L82: '2501' -> '2340.6'
$ws.Cells.Item(82, 13).Value = -1505.8334  # M82: '-1179.4' -> '-1505.8334'
$ws.Cells.Item(82, 14).Value = -3062.6  # N82: '-3223' -> '-3062.6'

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 2082.182  # H85: '1900.625' -> '2082.182'
$ws.Cells.Item(85, 9).Value = 1866.8334  # I85: '1540.4' -> '1866.8334'
$ws.Cells.Item(85, 10).Value = 2340.6  # J85: '2501' -> '2340.6'
$ws.Cells.Item(85, 11).Value = 1866.8334  # K85: '1540.4' -> '1866.8334'
$ws.Cells.Item(85, 12).Value = 2340.6  # L85: '2501' -> '2340.6'
$ws.Cells.Item(85, 13).Value = -618.8334  # M85: '-292.4000000000001' -> '-618.8334'
$ws.Cells.Item(85, 14).Value = -4836.6  # N85: '-4997' -> '-4836.6'

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 4799.8  # H93: '5000' -> '4799.8'
$ws.Cells.Item(93, 9).Value = 4799.8  # I93: '5000' -> '4799.8'
$ws.Cells.Item(93, 11).Value = 4799.8  # K93: '5000' -> '4799.8'
$ws.Cells.Item(93, 13).Value = -3551.8  # M93: '-3752' -> '-3551.8'

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 4034.3333  # H122: '4049.75' -> '4034.3333'
$ws.Cells.Item(122, 9).Value = 3904  # I122: '0' -> '3904'
$ws.Cells.Item(122, 10).Value = 4099.5  # J122: '4049.75' -> '4099.5'
$ws.Cells.Item(122, 11).Value = 11712  # K122: '0' -> '11712'
$ws.Cells.Item(122, 12).Value = 12298.5  # L122: '12149.25' -> '12298.5'
$ws.Cells.Item(122, 13).Value = -9262  # M122: None -> '-9262'
$ws.Cells.Item(122, 14).Value = -17198.5  # N122: '-17049.25' -> '-17198.5'

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 10000  # H132: '0' -> '10000'
$ws.Cells.Item(132, 9).Value = 10000  # I132: '0' -> '10000'
$ws.Cells.Item(132, 11).Value = 30000  # K132: '0' -> '30000'
$ws.Cells.Item(132, 13).Value = -27470  # M132: None -> '-27470'

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(139, 8).Value = 97977  # H139: '0' -> '97977'
$ws.Cells.Item(139, 10).Value = 97977  # J139: '0' -> '97977'
$ws.Cells.Item(139, 12).Value = 97977  # L139: '0' -> '97977'
$ws.Cells.Item(139, 14).Value = -108257  # N139: None -> '-108257'

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 20859  # H54: '15502.5' -> '20859'
$ws.Cells.Item(54, 9).Value = 16333.333  # I54: '15502.5' -> '16333.333'
$ws.Cells.Item(54, 10).Value = 25384.666  # J54: '0' -> '25384.666'
$ws.Cells.Item(54, 11).Value = 16333.333  # K54: '15502.5' -> '16333.333'
$ws.Cells.Item(54, 12).Value = 25384.666  # L54: '0' -> '25384.666'
$ws.Cells.Item(54, 13).Value = -15813.333  # M54: '-14982.5' -> '-15813.333'
$ws.Cells.Item(54, 14).Value = -26424.666  # N54: None -> '-26424.666'

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(95, 8).Value = 35999  # H95: '0' -> '35999'
$ws.Cells.Item(95, 10).Value = 35999  # J95: '0' -> '35999'
$ws.Cells.Item(95, 12).Value = 35999  # L95: '0' -> '35999'
$ws.Cells.Item(95, 14).Value = -41491  # N95: None -> '-41491'

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(104, 8).Value = 30000  # H104: '21500' -> '30000'
$ws.Cells.Item(104, 10).Value = 30000  # J104: '21500' -> '30000'
$ws.Cells.Item(104, 12).Value = 30000  # L104: '21500' -> '30000'
$ws.Cells.Item(104, 14).Value = -36988  # N104: '-28488' -> '-36988'

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 9).Value = 1000  # I122: '0' -> '1000'
$ws.Cells.Item(122, 11).Value = 3000  # K122: '0' -> '3000'
$ws.Cells.Item(122, 13).Value = -550  # M122: None -> '-550'

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 13625  # H132: '7751.8335' -> '13625'
$ws.Cells.Item(132, 9).Value = 13250  # I132: '6034' -> '13250'
$ws.Cells.Item(132, 10).Value = 14000  # J132: '11187.5' -> '14000'
$ws.Cells.Item(132, 11).Value = 39750  # K132: '18102' -> '39750'
$ws.Cells.Item(132, 12).Value = 42000  # L132: '33562.5' -> '42000'
$ws.Cells.Item(132, 13).Value = -37220  # M132: '-15572' -> '-37220'
$ws.Cells.Item(132, 14).Value = -47060  # N132: '-38622.5' -> '-47060'
